$wb = $excel.ActiveWorkbook

# --- Sheet "Forecast Comparison": update columns B (Prophet Forecast) and D (yhat_upper) ---
$ws1 = $wb.Worksheets.Item("Forecast Comparison")

$bValues = @{
    3  = 3
    6  = 4
    10 = 1
    11 = 4
    13 = 5
}

$dValues = @{
    2  = 17.64042441662114
    3  = 16.6201965439883
    4  = 16.6548992654382
    5  = 16.75117441384114
    6  = 18.22392699872321
    7  = 15.89800956019089
    8  = 12.93572156289842
    9  = 11.85865455729252
    10 = 14.31037386445029
    11 = 18.24644518752616
    12 = 18.08077902055489
    13 = 18.7823393435219
    14 = 15.78982458671447
    15 = 13.93146918752713
    16 = 11.23525802650791
    17 = 8.771486893698055
    18 = 5.785366465743787
    19 = 3.351765655209337
    20 = 3.492528385869525
    21 = 5.506021936574759
}

foreach ($row in $bValues.Keys) {
    $ws1.Cells.Item($row, 2).Value = $bValues[$row]
}

foreach ($row in $dValues.Keys) {
    $ws1.Cells.Item($row, 4).Value = $dValues[$row]
}

# --- Sheet "Summary": update Total Forecast (16 Weeks) and Total Forecast (4 Weeks) ---
# These cells store numeric-looking values as text, so force text format before writing.
$ws2 = $wb.Worksheets.Item("Summary")
$ws2.Range("B9").NumberFormat = "@"
$ws2.Range("B9").Value = "38"
$ws2.Range("B11").NumberFormat = "@"
$ws2.Range("B11").Value = "13"
